$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "27.890.72"
Set-TextValue "E2" "  -0.35%  "

# Row 3
Set-TextValue "D3" "1.632.16"
Set-TextValue "E3" "  -0.94%  "

# Row 4
Set-TextValue "E4" "  +0.03%  "

# Row 5
Set-TextValue "D5" "211.74"
Set-TextValue "E5" "  -0.88%  "

# Row 6
Set-TextValue "D6" "0.522"
Set-TextValue "E6" "  -0.99%  "

# Row 7
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.04%  "

# Row 8
Set-TextValue "D8" "23.23"
Set-TextValue "E8" "  -1.10%  "

# Row 9
Set-TextValue "D9" "0.257"
Set-TextValue "E9" "  -3.33%  "

# Row 10
Set-TextValue "E10" "  -0.26%  "

# Row 11
Set-TextValue "D11" "0.0879"
Set-TextValue "E11" "  +0.68%  "

# Row 12
Set-TextValue "D12" "1.863.89"
Set-TextValue "E12" "  -0.90%  "

# Row 13
Set-TextValue "D13" "1.642.69"
Set-TextValue "E13" "  -0.25%  "

# Row 14
Set-TextValue "E14" "  -0.53%  "

# Row 15
Set-TextValue "D15" "0.565"
Set-TextValue "E15" "  +0.20%  "

# Row 16
Set-TextValue "D16" "65.17"
Set-TextValue "E16" "  -0.71%  "

# Row 17
Set-TextValue "D17" "27.888.43"
Set-TextValue "E17" "  -0.31%  "

# Row 18
Set-TextValue "D18" "229.87"
Set-TextValue "E18" "  -1.39%  "

# Row 19
Set-TextValue "E19" "  -0.20%  "

# Row 20
Set-TextValue "E20" "  -2.45%  "

# Row 22
Set-TextValue "E22" "  -0.87%  "

# Row 24
Set-TextValue "E24" "  -4.11%  "

# Row 25
Set-TextValue "D25" "153.52"
Set-TextValue "E25" "  +0.57%  "

# Row 26
Set-TextValue "E26" "  +0.48%  "

# Row 27
Set-TextValue "E27" "  -0.75%  "

# Row 28
Set-TextValue "E28" "  -0.97%  "

# Row 29
Set-TextValue "E29" "  +0.04%  "

# Row 30
Set-TextValue "E30" "  -0.96%  "

# Row 31
Set-TextValue "E31" "  -0.30%  "

# Row 32
Set-TextValue "D32" "3.38"
Set-TextValue "E32" "  +0.56%  "

# Row 33
Set-TextValue "D33" "3.07"
Set-TextValue "E33" "  -1.03%  "

# Row 34
Set-TextValue "D34" "1.395.26"
Set-TextValue "E34" "  -3.75%  "

# Row 35
Set-TextValue "E35" "  -0.04%  "

# Row 36
Set-TextValue "D36" "1.01"
Set-TextValue "E36" "  +9.24%  "

# Row 37
Set-TextValue "E37" "  +1.45%  "

# Row 38
Set-TextValue "E38" "  +0.43%  "

# Row 39
Set-TextValue "D39" "0.559"
Set-TextValue "E39" "  -0.43%  "

# Row 40
Set-TextValue "D40" "0.870"
Set-TextValue "E40" "  -2.25%  "

# Row 41
Set-TextValue "E41" "  -0.22%  "

# Row 42
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.05%  "

# Row 43
Set-TextValue "D43" "66.86"
Set-TextValue "E43" "  -3.75%  "

# Row 44
Set-TextValue "E44" "  +3.02%  "

# Row 45
Set-TextValue "E45" "  +0.62%  "

# Row 46
Set-TextValue "E46" "  -1.12%  "

# Row 47
Set-TextValue "D47" "1.773.59"
Set-TextValue "E47" "  -0.90%  "

# Row 48
Set-TextValue "D48" "87.67"
Set-TextValue "E48" "  -1.50%  "

# Row 49
Set-TextValue "E49" "  -0.80%  "

# Row 50
Set-TextValue "D50" "0.0507"
Set-TextValue "E50" "  -0.12%  "

# Row 51
Set-TextValue "D51" "7.48"
Set-TextValue "E51" "  -3.15%  "
